$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column C (eps_growth) values for rows 28-122, refreshed FactSet pull ---
$ws.Cells.Item(28,3).Value2 = 15.0441
$ws.Cells.Item(29,3).Value2 = 10.1196
$ws.Cells.Item(30,3).Value2 = 5.0726500000000003
$ws.Cells.Item(31,3).Value2 = -12.729900000000001
$ws.Cells.Item(32,3).Value2 = -19.738499999999998
$ws.Cells.Item(33,3).Value2 = -18.8416
$ws.Cells.Item(34,3).Value2 = -20.462499999999999
$ws.Cells.Item(35,3).Value2 = -9.0701000000000001
$ws.Cells.Item(36,3).Value2 = -2.34884
$ws.Cells.Item(37,3).Value2 = 5.5220700000000003
$ws.Cells.Item(38,3).Value2 = 10.2902
$ws.Cells.Item(39,3).Value2 = 13.514799999999999
$ws.Cells.Item(40,3).Value2 = 20.1631
$ws.Cells.Item(41,3).Value2 = 20.9313
$ws.Cells.Item(42,3).Value2 = 25.1936
$ws.Cells.Item(43,3).Value2 = 28.482800000000001
$ws.Cells.Item(44,3).Value2 = 26.561199999999999
$ws.Cells.Item(45,3).Value2 = 20.5869
$ws.Cells.Item(46,3).Value2 = 21.6586
$ws.Cells.Item(47,3).Value2 = 17.5289
$ws.Cells.Item(48,3).Value2 = 16.545400000000001
$ws.Cells.Item(49,3).Value2 = 13.2509
$ws.Cells.Item(50,3).Value2 = 20.065300000000001
$ws.Cells.Item(51,3).Value2 = 18.9617
$ws.Cells.Item(52,3).Value2 = 20.684200000000001
$ws.Cells.Item(53,3).Value2 = 17.364100000000001
$ws.Cells.Item(54,3).Value2 = 9.9834099999999992
$ws.Cells.Item(55,3).Value2 = 7.9380300000000004
$ws.Cells.Item(56,3).Value2 = 1.53691
$ws.Cells.Item(57,3).Value2 = -0.50114000000000003
$ws.Cells.Item(58,3).Value2 = -7.1195500000000003
$ws.Cells.Item(59,3).Value2 = -12.5177
$ws.Cells.Item(60,3).Value2 = -17.616599999999998
$ws.Cells.Item(61,3).Value2 = -37.253599999999999
$ws.Cells.Item(62,3).Value2 = -44.868899999999996
$ws.Cells.Item(63,3).Value2 = -47.591700000000003
$ws.Cells.Item(64,3).Value2 = -43.181399999999996
$ws.Cells.Item(65,3).Value2 = 14.3802
$ws.Cells.Item(66,3).Value2 = 51.725099999999998
$ws.Cells.Item(67,3).Value2 = 83.839600000000004
$ws.Cells.Item(68,3).Value2 = 93.140500000000003
$ws.Cells.Item(69,3).Value2 = 36.615200000000002
$ws.Cells.Item(70,3).Value2 = 27.8963
$ws.Cells.Item(71,3).Value2 = 19.7849
$ws.Cells.Item(72,3).Value2 = 16.312799999999999
$ws.Cells.Item(73,3).Value2 = 16.4755
$ws.Cells.Item(74,3).Value2 = 11.7433
$ws.Cells.Item(75,3).Value2 = 10.3254
$ws.Cells.Item(76,3).Value2 = 6.8908500000000004
$ws.Cells.Item(77,3).Value2 = 1.9538800000000001
$ws.Cells.Item(78,3).Value2 = 1.64524
$ws.Cells.Item(79,3).Value2 = 1.54498
$ws.Cells.Item(80,3).Value2 = 4.2165299999999997
$ws.Cells.Item(81,3).Value2 = 8.0737500000000004
$ws.Cells.Item(82,3).Value2 = 9.1574299999999997
$ws.Cells.Item(83,3).Value2 = 10.597300000000001
$ws.Cells.Item(84,3).Value2 = 9.7029700000000005
$ws.Cells.Item(85,3).Value2 = 0.521424
$ws.Cells.Item(86,3).Value2 = -1.67926
$ws.Cells.Item(87,3).Value2 = -5.9399100000000002
$ws.Cells.Item(88,3).Value2 = -6.9796899999999997
$ws.Cells.Item(89,3).Value2 = -1.7497499999999999
$ws.Cells.Item(90,3).Value2 = -2.7810000000000001
$ws.Cells.Item(91,3).Value2 = -1.05887
$ws.Cells.Item(92,3).Value2 = -2.0750299999999999
$ws.Cells.Item(93,3).Value2 = -0.114522
$ws.Cells.Item(94,3).Value2 = 5.4486699999999999
$ws.Cells.Item(95,3).Value2 = 8.8250799999999998
$ws.Cells.Item(96,3).Value2 = 11.8857
$ws.Cells.Item(97,3).Value2 = 10.1342
$ws.Cells.Item(98,3).Value2 = 11.241099999999999
$ws.Cells.Item(99,3).Value2 = 14.6891
$ws.Cells.Item(100,3).Value2 = 19.074200000000001
$ws.Cells.Item(101,3).Value2 = 24.1663
$ws.Cells.Item(102,3).Value2 = 18.437799999999999
$ws.Cells.Item(103,3).Value2 = 11.825699999999999
$ws.Cells.Item(104,3).Value2 = 4.1870399999999997
$ws.Cells.Item(105,3).Value2 = -1.60236
$ws.Cells.Item(106,3).Value2 = -9.1237899999999996
$ws.Cells.Item(107,3).Value2 = -18.4575
$ws.Cells.Item(108,3).Value2 = -19.23
$ws.Cells.Item(109,3).Value2 = -19.804500000000001
$ws.Cells.Item(110,3).Value2 = 2.1103700000000001
$ws.Cells.Item(111,3).Value2 = 37.806699999999999
$ws.Cells.Item(112,3).Value2 = 55.825099999999999
$ws.Cells.Item(113,3).Value2 = 70.399100000000004
$ws.Cells.Item(114,3).Value2 = 52.158299999999997
$ws.Cells.Item(115,3).Value2 = 21.648599999999998
$ws.Cells.Item(116,3).Value2 = 8.2561
$ws.Cells.Item(117,3).Value2 = -2.3487300000000002
$ws.Cells.Item(118,3).Value2 = -8.3791399999999996
$ws.Cells.Item(119,3).Value2 = -3.31921
$ws.Cells.Item(120,3).Value2 = -0.67586900000000005
$ws.Cells.Item(121,3).Value2 = 2.1248300000000002
$ws.Cells.Item(122,3).Value2 = 7.8045
$ws.Cells.Item(123,3).Value2 = 7.7673699999999997

# --- Append new quarterly row 124 (date 12/31/2024) ---
$ws.Cells.Item(124, 1).Value2 = 45657
$ws.Cells.Item(124, 1).NumberFormat = "m/d/yy"
$ws.Cells.Item(124, 2).Value2 = 51
$ws.Cells.Item(124, 4).Value2 = 0

# --- Column A width grew slightly to fit the new data ---
$ws.Columns.Item(1).ColumnWidth = 9.71

# --- Restore the active selection used while reviewing the refreshed data ---
$ws.Range("C121").Select()
